# New crime data collected - weekly update for week of 11/20/2023-11/26/2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header strings: bump volume/report number and shift the reporting week.
# These live inside rich-text shared strings; editing the Characters range
# updates just the relevant run's text.
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "47"
$ws.Range("C9").Characters(27, 10).Text = "11/20/2023"
$ws.Range("C9").Characters(48, 10).Text = "11/26/2023"

# ---------------------------------------------------------------------------
# Cells that flip between the "0"/"***.*" text placeholders (style 14) and
# real numbers (style 15/16). Copying format+value from a donor cell that
# already carries the desired style keeps the shared-string / numFmt table
# clean (avoids spawning new styles), then we overwrite the value where the
# target is numeric.
# ---------------------------------------------------------------------------

# Row 15: C15 text("0") -> number 3
$ws.Range("F15").Copy($ws.Range("C15"))
$ws.Range("C15").Value2 = 3

# Row 18: C18 text("0") -> number 2
$ws.Range("F18").Copy($ws.Range("C18"))
$ws.Range("C18").Value2 = 2

# Row 22: D22 number(1) -> text("0"); E22 number(-100) -> text("***.*")
$ws.Range("D18").Copy($ws.Range("D22"))
$ws.Range("E18").Copy($ws.Range("E22"))

# Row 23: D23 number(1) -> text("0"); E23 number(0) -> text("***.*")
$ws.Range("D18").Copy($ws.Range("D23"))
$ws.Range("E18").Copy($ws.Range("E23"))

# Row 26: C26 text("0") -> number 3
$ws.Range("F26").Copy($ws.Range("C26"))
$ws.Range("C26").Value2 = 3

# Row 27: C27 text("0") -> number 1; F27 text("0") -> number 1
#         G27 number(1) -> text("0"); H27 number(-100) -> text("***.*")
$ws.Range("G26").Copy($ws.Range("C27"))
$ws.Range("C27").Value2 = 1
$ws.Range("G26").Copy($ws.Range("F27"))
$ws.Range("F27").Value2 = 1
$ws.Range("D18").Copy($ws.Range("G27"))
$ws.Range("E18").Copy($ws.Range("H27"))

# Row 28: C28 number(1) -> text("0")
$ws.Range("D18").Copy($ws.Range("C28"))

# Row 29: C29 number(1) -> text("0")
$ws.Range("D18").Copy($ws.Range("C29"))

# ---------------------------------------------------------------------------
# Plain numeric value updates (no style change required).
# ---------------------------------------------------------------------------

# Row 15
$ws.Range("F15").Value2 = 4
$ws.Range("H15").Value2 = 300
$ws.Range("I15").Value2 = 17
$ws.Range("K15").Value2 = 30.769230769230
$ws.Range("L15").Value2 = 21.428571428571
$ws.Range("M15").Value2 = 21.428571428571
$ws.Range("N15").Value2 = -22.727272727272

# Row 16
$ws.Range("C16").Value2 = 2
$ws.Range("D16").Value2 = 2
$ws.Range("E16").Value2 = 0
$ws.Range("F16").Value2 = 8
$ws.Range("G16").Value2 = 10
$ws.Range("H16").Value2 = -20
$ws.Range("I16").Value2 = 100
$ws.Range("J16").Value2 = 99
$ws.Range("K16").Value2 = 1.010101010101
$ws.Range("L16").Value2 = 88.679245283018
$ws.Range("M16").Value2 = -52.60663507109
$ws.Range("N16").Value2 = -84

# Row 17
$ws.Range("C17").Value2 = 9
$ws.Range("E17").Value2 = 50
$ws.Range("F17").Value2 = 16
$ws.Range("G17").Value2 = 22
$ws.Range("H17").Value2 = -27.272727272727
$ws.Range("I17").Value2 = 217
$ws.Range("J17").Value2 = 196
$ws.Range("K17").Value2 = 10.714285714285
$ws.Range("L17").Value2 = 13.020833333333
$ws.Range("M17").Value2 = 48.630136986301
$ws.Range("N17").Value2 = -13.2

# Row 18
$ws.Range("F18").Value2 = 4
$ws.Range("G18").Value2 = 2
$ws.Range("H18").Value2 = 100
$ws.Range("I18").Value2 = 49
$ws.Range("K18").Value2 = -41.666666666666
$ws.Range("L18").Value2 = -25.757575757575
$ws.Range("M18").Value2 = -79.148936170212
$ws.Range("N18").Value2 = -91.373239436619

# Row 19
$ws.Range("C19").Value2 = 4
$ws.Range("D19").Value2 = 3
$ws.Range("E19").Value2 = 33.333333333333
$ws.Range("I19").Value2 = 193
$ws.Range("J19").Value2 = 215
$ws.Range("K19").Value2 = -10.232558139534
$ws.Range("L19").Value2 = 21.383647798742
$ws.Range("M19").Value2 = -12.669683257918
$ws.Range("N19").Value2 = -47.554347826087

# Row 20
$ws.Range("C20").Value2 = 3
$ws.Range("D20").Value2 = 3
$ws.Range("E20").Value2 = 0
$ws.Range("F20").Value2 = 14
$ws.Range("G20").Value2 = 10
$ws.Range("H20").Value2 = 40
$ws.Range("I20").Value2 = 145
$ws.Range("J20").Value2 = 142
$ws.Range("K20").Value2 = 2.112676056338
$ws.Range("L20").Value2 = 57.608695652173
$ws.Range("M20").Value2 = 35.514018691588
$ws.Range("N20").Value2 = -92.654508611955

# Row 21 (TOTAL)
$ws.Range("C21").Value2 = 23
$ws.Range("D21").Value2 = 14
$ws.Range("E21").Value2 = 64.285714285714
$ws.Range("F21").Value2 = 60
$ws.Range("G21").Value2 = 62
$ws.Range("H21").Value2 = -3.225806451612
$ws.Range("I21").Value2 = 725
$ws.Range("J21").Value2 = 754
$ws.Range("K21").Value2 = -3.846153846153
$ws.Range("L21").Value2 = 23.931623931623
$ws.Range("M21").Value2 = -23.117709437963
$ws.Range("N21").Value2 = -80.996068152031

# Row 22
$ws.Range("L22").Value2 = 42.857142857142

# Row 23
$ws.Range("C23").Value2 = 3
$ws.Range("F23").Value2 = 8
$ws.Range("G23").Value2 = 6
$ws.Range("H23").Value2 = 33.333333333333
$ws.Range("I23").Value2 = 111
$ws.Range("K23").Value2 = 23.333333333333
$ws.Range("L23").Value2 = 35.365853658536
$ws.Range("M23").Value2 = 101.818181818182

# Row 24
$ws.Range("C24").Value2 = 5
$ws.Range("D24").Value2 = 3
$ws.Range("E24").Value2 = 66.666666666666
$ws.Range("F24").Value2 = 25
$ws.Range("G24").Value2 = 39
$ws.Range("H24").Value2 = -35.897435897435
$ws.Range("I24").Value2 = 489
$ws.Range("J24").Value2 = 501
$ws.Range("K24").Value2 = -2.395209580838
$ws.Range("L24").Value2 = 23.797468354430
$ws.Range("M24").Value2 = -0.609756097560

# Row 25
$ws.Range("D25").Value2 = 4
$ws.Range("E25").Value2 = 100
$ws.Range("F25").Value2 = 28
$ws.Range("G25").Value2 = 22
$ws.Range("H25").Value2 = 27.272727272727
$ws.Range("I25").Value2 = 332
$ws.Range("J25").Value2 = 343
$ws.Range("K25").Value2 = -3.206997084548
$ws.Range("L25").Value2 = 19.85559566787
$ws.Range("M25").Value2 = -20.952380952381

# Row 26
$ws.Range("F26").Value2 = 4
$ws.Range("H26").Value2 = 300
$ws.Range("I26").Value2 = 21
$ws.Range("K26").Value2 = 0
$ws.Range("L26").Value2 = -8.695652173913

# Row 27
$ws.Range("I27").Value2 = 25
$ws.Range("K27").Value2 = 8.695652173913
$ws.Range("L27").Value2 = -10.714285714285
